# Fill in the Preconditions (E), Method Inputs (F) and Expected Result (G)
# columns for the __init__ / __str__ / get_service_charges test cases on the
# ChequingAccount unit test plan. Values are entered column-by-column
# (E7:E14, then F7:F14, then G7:G14) to mirror how the sheet was authored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E: Preconditions ---------------------------------------------
$ws.Range("E7").Value = "None"
$ws.Range("E8").Value = "None"
$ws.Range("E9").Value = "None"
$ws.Range("E10").Value = "None"
$ws.Range("E11").Value = "Account is created with a balance greater than the overdraft limit"
$ws.Range("E12").Value = "Account is created with a balance less thatn the overdraft limit"
$ws.Range("E13").Value = "Account is created with a balance equal to the overdraft limit"
$ws.Range("E14").Value = "Account is created with specific values"

# --- Column F: Method Inputs ----------------------------------------------
$ws.Range("F7").Value = "account_number = 112233, client_number = 889900, balance = 500, date_created = date.today(), overdraft_limit = -100, overdraft_rate = 0.05"
$ws.Range("F8").Value = 'account_number = 112233, client_number = 889900, balance = 500, date_created = date.today(), overdraft_limit = "invalid", overdraft_rate = 0.05'
$ws.Range("F9").Value = 'account_number = 112233, client_number = 889900, balance = 500, date_created = date.today(), overdraft_limit = -100, overdraft_rate = "invalid"'
$ws.Range("F10").Value = 'account_number = 112233, client_number = 889900, balance = 500, date_created = "invalid date", overdraft_limit = -100, overdraft_rate = 0.05'
$ws.Range("F11").Value = "balance = 0, overdraft_limit = -100, overdraft_rate = 0.05"
$ws.Range("F12").Value = "balance = -600, overdraft_limit = -100, overdraft_rate = 0.05"
$ws.Range("F13").Value = "balance = -100, overdraft_limit = -100, overdraft_rate = 0.05"
$ws.Range("F14").Value = "account_number = 112233, client_number = 889900, balance = 1559.49, date_created=date.today(), overdraft_limit = -15.00, overdraft_rate = 0.05"

# --- Column G: Expected Result --------------------------------------------
$ws.Range("G7").Value = "Account initialized with correct attribute values."
$ws.Range("G8").Value = "Overdraft limit is set to default value which is -100 or raises a ValueError."
$ws.Range("G9").Value = "Overdraft rate is set to default value which is 0.05 or raises a ValueError."
$ws.Range("G10").Value = "Raises ValueError for invalid date type."
$ws.Range("G11").Value = "Total service charges calculated = `$0.50"
$ws.Range("G12").Value = "Total service charges calculated as `$0.50 + (amount over limit * rate)."
$ws.Range("G13").Value = "Total service charges calculated as `$0.50."
$ws.Range("G14").Value = 'Returns string: "Account Number: 112233 Balance: $1,559.49\nOverdraft Limit: $-15.00 Overdraft Rate: 5.00% Account Type: Chequing"'

# A few of the newly-filled cells pick up the "top-medium-border" look used
# by the header-adjacent row 7/8 style (rather than the thinner border used
# lower in the table) once they carry real content.
$xlEdgeTop = 8
$xlMedium = -4138
$ws.Range("F9").Borders.Item($xlEdgeTop).Weight = $xlMedium
$ws.Range("E10").Borders.Item($xlEdgeTop).Weight = $xlMedium
$ws.Range("F10").Borders.Item($xlEdgeTop).Weight = $xlMedium

# Match the view state captured in the saved workbook.
$ws.Application.ActiveWindow.Zoom = 84
$ws.Range("G14").Select()
